# Update size in ids columns
#
# The placeholder IDs in the first ("ID") column of the Creature / Effect /
# Trait / Mechanic tables gain a separating hyphen (e.g. "CREXXX" becomes
# "CRE-XXX"). That one extra character is what makes column A on each of
# those sheets a bit wider (the column has best-fit/auto-size behaviour).

$wb = $excel.ActiveWorkbook

$wsCreatures = $wb.Worksheets.Item("Creatures")
$wsEffects   = $wb.Worksheets.Item("Effects")
$wsTraits    = $wb.Worksheets.Item("Traits")
$wsColors    = $wb.Worksheets.Item("Colors Overview")

# 1. Add the missing hyphen to the example/placeholder IDs.
$wsCreatures.Range("A2").Value = "CRE-XXX"
$wsEffects.Range("A2").Value   = "EFF-XXX"
$wsTraits.Range("A2").Value    = "TRA-XXX"
$wsColors.Range("A2").Value    = "MEC-XXX"

# 2. Re-fit column A on each sheet now that its content is one character
#    wider.
$wsCreatures.Columns.Item(1).ColumnWidth = 7.666666666666667
$wsEffects.Columns.Item(1).ColumnWidth   = 7.333333333333333
$wsTraits.Columns.Item(1).ColumnWidth    = 7.666666666666667
$wsColors.Columns.Item(1).ColumnWidth    = 8.166666666666666

# 3. The Effects sheet's active cell moved to A2 (the cell that was just
#    edited).
$wsEffects.Range("A2").Select()
